# Regenerate handback-status report:
#   - update the existing row (UUID 127c6d9f... -> 170a976a..., refreshed timestamps)
#   - append a brand-new row for a second file (ffff32cbb600...)
# across the Overview, zh-cn and de-de sheets, keeping the tables/ranges in sync.

$wb = $excel.ActiveWorkbook

$oldUuid = "127c6d9f-82d3-41eb-88e8-925348b722ca"
$newUuid = "170a976a-1720-47b4-8170-4c357abdc0e5"
$addUuid = "ffff32cbb600-8252-401b-afc2-e962e4a40b74"

$newSha = "f05706c3fed05b3e509b9b2cb9dc4f6afa05ae66"

# Helper: write a literal text value into a cell without Excel's automatic
# "True"/"False" -> Boolean coercion (round-trips through a text formula and
# a values-only paste so the result lands back in the sheet as shared text).
function Set-LiteralText($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Row 2: refresh the UUID / hyperlink / generate-date for the existing file.
$wsOv.Range("A2").Value = "$newUuid.md"

$wsOv.Range("B2").Hyperlinks.Delete()
$wsOv.Range("B2").Value = "e2e\$newUuid.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/$newSha/e2e/$newUuid.md", $null, $null, "e2e\$newUuid.md") | Out-Null

$wsOv.Range("G2").Value = "2016-08-13 23:27:23"
$wsOv.Range("G2").NumberFormat = $dateFormat

# Row 3: brand-new file.
$wsOv.Range("A3").Value = "$addUuid.md"

$wsOv.Range("B3").Value = "e2e\$addUuid.md"
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/$newSha/e2e/$addUuid.md", $null, $null, "e2e\$addUuid.md") | Out-Null

$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value = "Handed back: in sync with en-US"
$wsOv.Range("G3").Value = "2016-08-13 23:27:23"
$wsOv.Range("G3").NumberFormat = $dateFormat

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$newUuid.$newSha.zh-cn.xlf"

# Row 2: refresh UUID / xlf name / timestamps; hyperlinks keep their rIds.
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A2").Value = "$newUuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/$newSha/e2e/$newUuid.md", $null, $null, "$newUuid.md") | Out-Null

$wsZh.Range("G2").Value = $zhXlf
$wsZh.Range("H2").Value = "2016-08-13 23:27:15"
$wsZh.Range("H2").NumberFormat = $dateFormat

$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Range("I2").Value = "$newUuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8a5f1f2d061862bd0e637d07c72026e16019e71a/e2e/$newUuid.md", $null, $null, "$newUuid.md") | Out-Null

$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = "2016-08-13 23:27:43"
$wsZh.Range("K2").NumberFormat = $dateFormat

# Row 3: brand-new file.
$wsZh.Range("A3").Value = "$addUuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/$newSha/e2e/$addUuid.md", $null, $null, "$addUuid.md") | Out-Null

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
Set-LiteralText $wsZh.Range("F3") "True"

$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = "2016-08-13 23:27:15"
$wsZh.Range("H3").NumberFormat = $dateFormat

$wsZh.Range("I3").Value = "$addUuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8a5f1f2d061862bd0e637d07c72026e16019e71a/e2e/$addUuid.md", $null, $null, "$addUuid.md") | Out-Null

$wsZh.Range("J3").Value = $zhXlf
$wsZh.Range("K3").Value = "2016-08-13 23:27:43"
$wsZh.Range("K3").NumberFormat = $dateFormat

Set-LiteralText $wsZh.Range("M3") "True"
Set-LiteralText $wsZh.Range("O3") "False"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlf = "$newUuid.$newSha.de-de.xlf"

# Row 2: refresh UUID / xlf name / timestamps; hyperlinks keep their rIds.
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A2").Value = "$newUuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/$newSha/e2e/$newUuid.md", $null, $null, "$newUuid.md") | Out-Null

$wsDe.Range("G2").Value = $deXlf
$wsDe.Range("H2").Value = "2016-08-13 23:27:23"
$wsDe.Range("H2").NumberFormat = $dateFormat

$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Range("I2").Value = "$newUuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/240ef9aec38f0dcc70deb4e7ad048ee61c7c89f4/e2e/$newUuid.md", $null, $null, "$newUuid.md") | Out-Null

$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = "2016-08-13 23:27:53"
$wsDe.Range("K2").NumberFormat = $dateFormat

# Row 3: brand-new file.
$wsDe.Range("A3").Value = "$addUuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/$newSha/e2e/$addUuid.md", $null, $null, "$addUuid.md") | Out-Null

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
Set-LiteralText $wsDe.Range("F3") "True"

$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = "2016-08-13 23:27:23"
$wsDe.Range("H3").NumberFormat = $dateFormat

$wsDe.Range("I3").Value = "$addUuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/240ef9aec38f0dcc70deb4e7ad048ee61c7c89f4/e2e/$addUuid.md", $null, $null, "$addUuid.md") | Out-Null

$wsDe.Range("J3").Value = $deXlf
$wsDe.Range("K3").Value = "2016-08-13 23:27:53"
$wsDe.Range("K3").NumberFormat = $dateFormat

Set-LiteralText $wsDe.Range("M3") "True"
Set-LiteralText $wsDe.Range("O3") "False"

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))

"Handback status report regenerated."
